# Adds newly-tracked movies/shows to the "media" (sheet1) and "gay" (sheet2)
# worksheets, appending rows to the bottom of each table.
#
# The order in which the brand-new text values below are written matters:
# it determines the order new entries are appended to the shared string
# table, which mirrors how the rows were originally typed into the two
# sheets (interleaved between "media" and "gay").

$wb = $excel.ActiveWorkbook
$media = $wb.Worksheets.Item(1)   # "media" sheet (sheet1.xml)
$gay   = $wb.Worksheets.Item(2)   # "gay" sheet (sheet2.xml)

# -- gay!A139:B139 --------------------------------------------------------
$gay.Range("A139").Value = "Maurice"
$gay.Range("B139").Value = "tt0093512"

# -- media!A557:B557 -------------------------------------------------------
$media.Range("A557").Value = "Evolution"
$media.Range("B557").Value = "tt4291590"

# -- media!A558:B558 -------------------------------------------------------
$media.Range("A558").Value = "Polar"
$media.Range("B558").Value = "tt4139588"

# -- media!A559:B559 -------------------------------------------------------
$media.Range("A559").Value = "A Star is Born"
$media.Range("B559").Value = "tt1517451"

# -- gay!A140:B140 ----------------------------------------------------------
$gay.Range("A140").Value = "Fair Haven"
$gay.Range("B140").Value = "tt3520216"

# -- gay!A141:B141 ----------------------------------------------------------
$gay.Range("A141").Value = "Gewoon Vrienden"
$gay.Range("B141").Value = "tt7901640"

# -- media!A560:B560 -------------------------------------------------------
$media.Range("A560").Value = "Mortal Engines"
$media.Range("B560").Value = "tt1571234"

# -- gay!B142, A143, A142, B143 (entered out of row order) -----------------
$gay.Range("B142").Value = "tt7008872"
$gay.Range("A143").Value = "Parting Glances"
$gay.Range("A142").Value = "Boy Erased"
$gay.Range("B143").Value = "tt0091725"

# -- media!A561:B561 -------------------------------------------------------
$media.Range("A561").Value = "Fantastic Beasts The Crimes of Grindelwald"
$media.Range("B561").Value = "tt4123430"

# -- media!A562:B562 -------------------------------------------------------
$media.Range("A562").Value = "The Boy Who Could Fly"
$media.Range("B562").Value = "tt0090768"

# -- media!A563:B563 -------------------------------------------------------
$media.Range("A563").Value = "The Last Mimzy"
$media.Range("B563").Value = "tt0768212"

# -- Fill in the "Type" column for every new row (all are Movies) ----------
$media.Range("C557").Value = "Movie"
$media.Range("C558").Value = "Movie"
$media.Range("C559").Value = "Movie"
$media.Range("C560").Value = "Movie"
$media.Range("C561").Value = "Movie"
$media.Range("C562").Value = "Movie"
$media.Range("C563").Value = "Movie"

$gay.Range("C139").Value = "Movie"
$gay.Range("C140").Value = "Movie"
$gay.Range("C141").Value = "Movie"
$gay.Range("C142").Value = "Movie"
$gay.Range("C143").Value = "Movie"

# -- Restore the view state (active cell / scroll position) ----------------
[void]$media.Activate()
[void]$media.Range("C560:C563").Select()

[void]$gay.Activate()
[void]$gay.Range("C142:C143").Select()

[void]$media.Activate()
